$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E5").Value = "12,5%"
$ws.Range("F5").Value = "87,5%"

$ws.Range("D7").Value = "12,5%"
$ws.Range("E7").Value = "12,5%"

$ws.Range("D11").Value = "12,5%"
$ws.Range("E11").Value = "12,5%"

$ws.Range("D12").Value = "12,5%"
$ws.Range("E12").Value = "12,5%"

$ws.Range("D13").Value = "12,5%"
$ws.Range("E13").Value = "12,5%"

$ws.Range("D14").Value = "12,5%"
$ws.Range("E14").Value = "12,5%"
